$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new matches (TURKEY - 1. LIG, 09:00, ROUND 29) were added to the sheet,
# inserted right before the existing "Rizespor vs Manisa FK" match (which was
# row 5 and is now pushed down to row 7). The new rows only have the basic
# match info (Date/Time/League/Home/Away/Round) and no odds columns.

# Insert two blank rows at row 5 - this pushes the old row 5 (and its odds
# data) down to row 7, exactly as in the target workbook.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).Insert()

# New row 5: Genclerbirligi vs Adanaspor AS
$ws.Range("A5").Value = "29/03/2023"
$ws.Range("B5").Value = "09:00"
$ws.Range("C5").Value = "TURKEY - 1. LIG"
$ws.Range("D5").Value = "Genclerbirligi"
$ws.Range("E5").Value = "Adanaspor AS"
$ws.Range("F5").Value = "ROUND 29"

# New row 6: Goztepe vs Yeni Malatyaspor
$ws.Range("A6").Value = "29/03/2023"
$ws.Range("B6").Value = "09:00"
$ws.Range("C6").Value = "TURKEY - 1. LIG"
$ws.Range("D6").Value = "Goztepe"
$ws.Range("E6").Value = "Yeni Malatyaspor"
$ws.Range("F6").Value = "ROUND 29"
